# "importamos e exportamos, só falta o histórico"
#
# The "status" column (F) stored boolean-ish text values "True"/"False".
# Relabel them to the Portuguese "ativo"/"inativo", and leave the
# selection on F2 (the top of that status column) instead of the
# previous M13 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data so we only touch the used range.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Column F holds the "status" values that used to be the literal
# strings "True"/"False"; rename them to "ativo"/"inativo" everywhere
# they occur, matching the whole cell contents only.
$statusRange = $ws.Range("F1:F$lastRow")
$statusRange.Replace("True", "ativo", 1)
$statusRange.Replace("False", "inativo", 1)

# Move the active selection to F2.
$ws.Range("F2").Select()
